# In Progress test_Employee_details test case
# Adds a new "Employee_Details" worksheet populated with employee detail
# test-data columns, and refreshes the "Sweta Arora" test-data row on the
# PIM_Add_Employee sheet (previously "Nishchay Angra").

$wb = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item("Login")
$pim = $wb.Worksheets.Item("PIM_Add_Employee")

# ---------------------------------------------------------------------
# 1) PIM_Add_Employee: refresh the sample employee row with new values
# ---------------------------------------------------------------------
$pim.Range("A2").Value = "Sweta"
$pim.Range("C2").Value = "Arora"
$pim.Range("D2").Value = 2965
$pim.Range("E2").Value = "Sweta@Arora_10"
$pim.Range("F2").Value = "Sweta@20"
$pim.Range("G2").Value = "Sweta@20"

# createusername (E2) now also becomes a mailto hyperlink, matching the
# existing createpassword (F2) hyperlink pattern.
$pim.Hyperlinks.Add($pim.Range("E2"), "mailto:Sweta@Arora_10") | Out-Null

$pim.Range("D2").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Add the new Employee_Details worksheet after PIM_Add_Employee
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Employee_Details"

$headers = @("otherid","drivinglicence","licexpmonth","nationality","licexpyear","licexpdate","maritalstatus","dobmonth","dobyear","dobdate","gender","bloodtype","testfield")
$cols    = @("A","B","C","D","E","F","G","H","I","J","K","L","M")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $addr = $cols[$i] + "1"
    $ws.Range($addr).Value = $headers[$i]
}

# Header formatting: reuse the exact existing styles from the workbook so
# the shared style table stays as close as possible to a hand-edited file.
# "bold + full border, vertically centered" (createusername/password/confirm headers)
foreach ($c in @("A","B","C","D","G","H")) {
    $pim.Range("E1").Copy() | Out-Null
    $ws.Range($c + "1").PasteSpecial(-4122) | Out-Null
}
# "bold + full border" (Login url/username/password headers)
foreach ($c in @("E","F","I","J")) {
    $login.Range("A1").Copy() | Out-Null
    $ws.Range($c + "1").PasteSpecial(-4122) | Out-Null
}
# "bold + left/right border only" (Login expected_title header)
foreach ($c in @("K","L")) {
    $login.Range("D1").Copy() | Out-Null
    $ws.Range($c + "1").PasteSpecial(-4122) | Out-Null
}
# "bold, vertically centered, no border" -- new style
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").VerticalAlignment = -4108

# Restore the header values (PasteSpecial formats-only shouldn't disturb
# them, but set again defensively in case of engine quirks).
for ($i = 0; $i -lt $headers.Length; $i++) {
    $addr = $cols[$i] + "1"
    $ws.Range($addr).Value = $headers[$i]
}

# Row 2 sample data
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "DL29AC2389"
$ws.Range("C2").Value = "August"
$ws.Range("D2").Value = "Indian"
$ws.Range("E2").Value = 2023
$ws.Range("F2").Value = 25
$ws.Range("G2").Value = "Single"
$ws.Range("H2").Value = "October"
$ws.Range("I2").Value = 1995
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = "Female"
$ws.Range("L2").Value = "B+"
$ws.Range("M2").Value = "Test"

# Data formatting: plain + full border for A2:J2 (reuse existing style)
foreach ($c in @("A","B","D","E","F","G","H","I","J")) {
    $pim.Range("A2").Copy() | Out-Null
    $ws.Range($c + "2").PasteSpecial(-4122) | Out-Null
}
# C2 (licexpmonth value) carries a date number format plus the full border
$pim.Range("A2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").NumberFormat = "mm-dd-yy"

# Reassign values once more after the format paste (PasteSpecial with
# xlPasteFormats should not touch values, but keep this deterministic).
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "DL29AC2389"
$ws.Range("C2").Value = "August"
$ws.Range("D2").Value = "Indian"
$ws.Range("E2").Value = 2023
$ws.Range("F2").Value = 25
$ws.Range("G2").Value = "Single"
$ws.Range("H2").Value = "October"
$ws.Range("I2").Value = 1995
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = "Female"
$ws.Range("L2").Value = "B+"
$ws.Range("M2").Value = "Test"

# Column widths (best-fit approximations of the authored widths)
$widths = @{ "B"=16.6640625; "C"=12.109375; "D"=13.33203125; "E"=9.44140625; "F"=9.5546875; "G"=11.88671875; "H"=11.44140625; "I"=9.44140625; "J"=9.5546875; "M"=11 }
foreach ($c in $widths.Keys) {
    $ws.Columns.Item($c).ColumnWidth = $widths[$c] - 0.8333333333
}

$ws.PageSetup.Orientation = 1

$ws.Range("M2").Select() | Out-Null

# Keep PIM_Add_Employee as the active/visible tab, matching the source
# workbook's tabSelected state.
$pim.Activate() | Out-Null
